$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.6828303333333333
$ws.Range("H2").Value = 2.048491
$ws.Range("I2").Value = 0.005533692967597834
$ws.Range("J2").Value = 0.005533692967597833
$ws.Range("M2").Value = 0.6828303333333333
$ws.Range("N2").Value = 2.048491
$ws.Range("O2").Value = 0.005533692967597834
$ws.Range("P2").Value = 0.005533692967597833
$ws.Range("Q2").Value = 0.4662572641201111
$ws.Range("R2").Value = 4.196315377081
$ws.Range("S2").Value = 0.00003062175785964173
$ws.Range("T2").Value = 0.00003062175785964171
$ws.Range("G3").Value = 0.6828303333333333
$ws.Range("H3").Value = 2.048491
$ws.Range("I3").Value = 0.005533692967597834
$ws.Range("J3").Value = 0.005533692967597833
$ws.Range("M3").Value = 49.58860766666667
$ws.Range("N3").Value = 148.765823
$ws.Range("O3").Value = 0.40186868702572
$ws.Range("P3").Value = 0.40186868702572
$ws.Range("Q3").Value = 33.86060550256589
$ws.Range("R3").Value = 304.745449523093
$ws.Range("S3").Value = 0.002223817927292002
$ws.Range("T3").Value = 0.002223817927292001
$ws.Range("G4").Value = 0.6828303333333333
$ws.Range("H4").Value = 2.048491
$ws.Range("I4").Value = 0.005533692967597834
$ws.Range("J4").Value = 0.005533692967597833
$ws.Range("M4").Value = 16.415838
$ws.Range("N4").Value = 49.247514
$ws.Range("O4").Value = 0.1330348153316153
$ws.Range("P4").Value = 0.1330348153316153
$ws.Range("Q4").Value = 11.209232133486
$ws.Range("R4").Value = 100.883089201374
$ws.Range("S4").Value = 0.0007361738220462359
$ws.Range("T4").Value = 0.0007361738220462358
$ws.Range("G5").Value = 0.6828303333333333
$ws.Range("H5").Value = 2.048491
$ws.Range("I5").Value = 0.005533692967597834
$ws.Range("J5").Value = 0.005533692967597833
$ws.Range("M5").Value = 56.70777633333334
$ws.Range("N5").Value = 170.123329
$ws.Range("O5").Value = 0.459562804675067
$ws.Range("P5").Value = 0.459562804675067
$ws.Range("Q5").Value = 38.72178981628211
$ws.Range("R5").Value = 348.496108346539
$ws.Range("S5").Value = 0.002543079460399955
$ws.Range("T5").Value = 0.002543079460399955
$ws.Range("G6").Value = 49.58860766666667
$ws.Range("H6").Value = 148.765823
$ws.Range("I6").Value = 0.40186868702572
$ws.Range("J6").Value = 0.40186868702572
$ws.Range("M6").Value = 0.6828303333333333
$ws.Range("N6").Value = 2.048491
$ws.Range("O6").Value = 0.005533692967597834
$ws.Range("P6").Value = 0.005533692967597833
$ws.Range("Q6").Value = 33.86060550256589
$ws.Range("R6").Value = 304.745449523093
$ws.Range("S6").Value = 0.002223817927292002
$ws.Range("T6").Value = 0.002223817927292001
$ws.Range("G7").Value = 49.58860766666667
$ws.Range("H7").Value = 148.765823
$ws.Range("I7").Value = 0.40186868702572
$ws.Range("J7").Value = 0.40186868702572
$ws.Range("M7").Value = 49.58860766666667
$ws.Range("N7").Value = 148.765823
$ws.Range("O7").Value = 0.40186868702572
$ws.Range("P7").Value = 0.40186868702572
$ws.Range("Q7").Value = 2459.030010318592
$ws.Range("R7").Value = 22131.27009286733
$ws.Range("S7").Value = 0.1614984416117761
$ws.Range("T7").Value = 0.1614984416117761
$ws.Range("G8").Value = 49.58860766666667
$ws.Range("H8").Value = 148.765823
$ws.Range("I8").Value = 0.40186868702572
$ws.Range("J8").Value = 0.40186868702572
$ws.Range("M8").Value = 16.415838
$ws.Range("N8").Value = 49.247514
$ws.Range("O8").Value = 0.1330348153316153
$ws.Range("P8").Value = 0.1330348153316153
$ws.Range("Q8").Value = 814.038550101558
$ws.Range("R8").Value = 7326.346950914023
$ws.Range("S8").Value = 0.05346252656602535
$ws.Range("T8").Value = 0.05346252656602535
$ws.Range("G9").Value = 49.58860766666667
$ws.Range("H9").Value = 148.765823
$ws.Range("I9").Value = 0.40186868702572
$ws.Range("J9").Value = 0.40186868702572
$ws.Range("M9").Value = 56.70777633333334
$ws.Range("N9").Value = 170.123329
$ws.Range("O9").Value = 0.459562804675067
$ws.Range("P9").Value = 0.459562804675067
$ws.Range("Q9").Value = 2812.059672242752
$ws.Range("R9").Value = 25308.53705018477
$ws.Range("S9").Value = 0.1846839009206266
$ws.Range("T9").Value = 0.1846839009206266
$ws.Range("G10").Value = 16.415838
$ws.Range("H10").Value = 49.247514
$ws.Range("I10").Value = 0.1330348153316153
$ws.Range("J10").Value = 0.1330348153316153
$ws.Range("M10").Value = 0.6828303333333333
$ws.Range("N10").Value = 2.048491
$ws.Range("O10").Value = 0.005533692967597834
$ws.Range("P10").Value = 0.005533692967597833
$ws.Range("Q10").Value = 11.209232133486
$ws.Range("R10").Value = 100.883089201374
$ws.Range("S10").Value = 0.0007361738220462359
$ws.Range("T10").Value = 0.0007361738220462358
$ws.Range("G11").Value = 16.415838
$ws.Range("H11").Value = 49.247514
$ws.Range("I11").Value = 0.1330348153316153
$ws.Range("J11").Value = 0.1330348153316153
$ws.Range("M11").Value = 49.58860766666667
$ws.Range("N11").Value = 148.765823
$ws.Range("O11").Value = 0.40186868702572
$ws.Range("P11").Value = 0.40186868702572
$ws.Range("Q11").Value = 814.038550101558
$ws.Range("R11").Value = 7326.346950914023
$ws.Range("S11").Value = 0.05346252656602535
$ws.Range("T11").Value = 0.05346252656602535
$ws.Range("G12").Value = 16.415838
$ws.Range("H12").Value = 49.247514
$ws.Range("I12").Value = 0.1330348153316153
$ws.Range("J12").Value = 0.1330348153316153
$ws.Range("M12").Value = 16.415838
$ws.Range("N12").Value = 49.247514
$ws.Range("O12").Value = 0.1330348153316153
$ws.Range("P12").Value = 0.1330348153316153
$ws.Range("Q12").Value = 269.479737242244
$ws.Range("R12").Value = 2425.317635180196
$ws.Range("S12").Value = 0.01769826209031698
$ws.Range("T12").Value = 0.01769826209031698
$ws.Range("G13").Value = 16.415838
$ws.Range("H13").Value = 49.247514
$ws.Range("I13").Value = 0.1330348153316153
$ws.Range("J13").Value = 0.1330348153316153
$ws.Range("M13").Value = 56.70777633333334
$ws.Range("N13").Value = 170.123329
$ws.Range("O13").Value = 0.459562804675067
$ws.Range("P13").Value = 0.459562804675067
$ws.Range("Q13").Value = 930.905669628234
$ws.Range("R13").Value = 8378.151026654106
$ws.Range("S13").Value = 0.06113785285322671
$ws.Range("T13").Value = 0.06113785285322671
$ws.Range("G14").Value = 56.70777633333334
$ws.Range("H14").Value = 170.123329
$ws.Range("I14").Value = 0.459562804675067
$ws.Range("J14").Value = 0.459562804675067
$ws.Range("M14").Value = 0.6828303333333333
$ws.Range("N14").Value = 2.048491
$ws.Range("O14").Value = 0.005533692967597834
$ws.Range("P14").Value = 0.005533692967597833
$ws.Range("Q14").Value = 38.72178981628211
$ws.Range("R14").Value = 348.496108346539
$ws.Range("S14").Value = 0.002543079460399955
$ws.Range("T14").Value = 0.002543079460399955
$ws.Range("G15").Value = 56.70777633333334
$ws.Range("H15").Value = 170.123329
$ws.Range("I15").Value = 0.459562804675067
$ws.Range("J15").Value = 0.459562804675067
$ws.Range("M15").Value = 49.58860766666667
$ws.Range("N15").Value = 148.765823
$ws.Range("O15").Value = 0.40186868702572
$ws.Range("P15").Value = 0.40186868702572
$ws.Range("Q15").Value = 2812.059672242752
$ws.Range("R15").Value = 25308.53705018477
$ws.Range("S15").Value = 0.1846839009206266
$ws.Range("T15").Value = 0.1846839009206266
$ws.Range("G16").Value = 56.70777633333334
$ws.Range("H16").Value = 170.123329
$ws.Range("I16").Value = 0.459562804675067
$ws.Range("J16").Value = 0.459562804675067
$ws.Range("M16").Value = 16.415838
$ws.Range("N16").Value = 49.247514
$ws.Range("O16").Value = 0.1330348153316153
$ws.Range("P16").Value = 0.1330348153316153
$ws.Range("Q16").Value = 930.905669628234
$ws.Range("R16").Value = 8378.151026654106
$ws.Range("S16").Value = 0.06113785285322671
$ws.Range("T16").Value = 0.06113785285322671
$ws.Range("G17").Value = 56.70777633333334
$ws.Range("H17").Value = 170.123329
$ws.Range("I17").Value = 0.459562804675067
$ws.Range("J17").Value = 0.459562804675067
$ws.Range("M17").Value = 56.70777633333334
$ws.Range("N17").Value = 170.123329
$ws.Range("O17").Value = 0.459562804675067
$ws.Range("P17").Value = 0.459562804675067
$ws.Range("Q17").Value = 3215.771896671361
$ws.Range("R17").Value = 28941.94707004225
$ws.Range("S17").Value = 0.2111979714408138
$ws.Range("T17").Value = 0.2111979714408138
